$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update Correspond Handoff/Handback Datetime for the first data row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 07:06:16"
$wsZhCn.Range("H2").Value = "2016-03-13 07:06:34"

# "de-de" sheet: update Correspond Handoff/Handback Datetime for the first data row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 07:06:20"
$wsDeDe.Range("H2").Value = "2016-03-13 07:06:41"
